# Apply the two substantive OOXML changes captured by the commit:
#
# 1. The table on slide 5 switches from the deck's custom "Table_0" table
#    style ({F8DA5D5B-D5D8-42C8-9E69-2DF968DDD4CF}) to the built-in
#    "Medium Style 2 - Accent 1" table style
#    ({0845A67C-2787-4A5E-A8FC-EE709F202FFD}).
#
# 2. The presentation's colour theme is swapped from the custom "Integral /
#    Red Violet" palette to the stock Office palette (fonts and effects are
#    identical between the two themes, only the 12 theme colours differ).

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 ---------------------------------------
$slide = $p.Slides.Item(5)
$tableShape = $slide.Shapes.Item(2)
$table = $tableShape.Table
$table.ApplyStyle("{0845A67C-2787-4A5E-A8FC-EE709F202FFD}")

# --- 2. Swap the theme colour scheme ----------------------------------
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink -> packed VBA RGB() values
# for the stock "Office" palette.
$officeThemeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

$colorScheme = $p.SlideMaster.Theme.ThemeColorScheme
for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}
